{"js": "// Fixed #418 Empty AQL expressions generate empty lines.\n// Remove the empty paragraph (with the orange-colored empty run) that\n// immediately follows the \"... :\" paragraph, so the empty AQL expression\n// no longer leaves a blank line behind.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text === \"\") {\n    paragraph.font.load(\"color\");\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text === \"\" && paragraph.font.color === \"#E36C0A\") {\n    paragraph.delete();\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Fixed #418 Empty AQL expressions generate empty lines.\n# Remove the empty paragraph (with the orange-colored empty run) that\n# immediately follows the \"... :\" paragraph, so the empty AQL expression\n# no longer leaves a blank line behind.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    $trimmed = $text.TrimEnd([char]13, [char]7)\n    if ($trimmed -eq \"\" -and $p.Range.Font.Color -eq -654262273) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
